$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 548.6514879859118
$ws.Range("D2").Value = 134.4137704555584
$ws.Range("F2").Value = 449
$ws.Range("G2").Value = 503
$ws.Range("H2").Value = 615
$ws.Range("C3").Value = 40.53964289067322
$ws.Range("D3").Value = 4.715725721416273
$ws.Range("F3").Value = 37.65
$ws.Range("G3").Value = 39.85
$ws.Range("H3").Value = 43.1
$ws.Range("C4").Value = 1.402195413855055
$ws.Range("D4").Value = 2.218147992283134
$ws.Range("G4").Value = 1.01
$ws.Range("H4").Value = 1.8
$ws.Range("C5").Value = 323.8575897991681
$ws.Range("D5").Value = 10.27855218193386
$ws.Range("F5").Value = 318
$ws.Range("G5").Value = 325.63
$ws.Range("H5").Value = 332.36
$ws.Range("C6").Value = 20.7228952216187
$ws.Range("D6").Value = 2.393059563722234
$ws.Range("E6").Value = 13.74
$ws.Range("F6").Value = 19.34
$ws.Range("G6").Value = 20.45
$ws.Range("H6").Value = 22.12
$ws.Range("C7").Value = -76.45389745753982
$ws.Range("D7").Value = 22.48866616808437
$ws.Range("C8").Value = 7.707435274623196
$ws.Range("D8").Value = 6.877757018260256
$ws.Range("C9").Value = 9.32302998759442
$ws.Range("D9").Value = 1.688001983131599
$ws.Range("C10").Value = 867.8303149157719
$ws.Range("D10").Value = 0.4611064813043405
$ws.Range("C11").Value = 0.5569493038346092
$ws.Range("D11").Value = 0.5905209443196567
$ws.Range("C12").Value = 22.68986415844869
$ws.Range("D12").Value = 12.27851648918861
$ws.Range("C13").Value = 0.6713688157821577
$ws.Range("D13").Value = 0.7481981734396599
$ws.Range("C14").Value = 1.826429444519897
$ws.Range("D14").Value = 1.666535742865267
$ws.Range("C15").Value = 93.85389745753963
$ws.Range("D15").Value = 22.48866616808437
$ws.Range("C16").Value = -85.73525626455496
$ws.Range("D16").Value = 20.24087520198974
$ws.Range("F16").Value = -101.4139268515822
$ws.Range("G16").Value = -85.41392685158225
$ws.Range("H16").Value = -68.14699179957641
$ws.Range("C17").Value = -78.02782098993175
$ws.Range("D17").Value = 24.89145279269826
$ws.Range("F17").Value = -92.39612087980606
$ws.Range("G17").Value = -75.2778545523916
$ws.Range("H17").Value = -57.3175485570292
